$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H259").Value = 6086
$ws.Range("H273").Value = 27129
$ws.Range("I273").Value = 1366
$ws.Range("H278").Value = 30006
$ws.Range("H287").Value = 57649
$ws.Range("H300").Value = 70463
$ws.Range("H301").Value = 69574
$ws.Range("H304").Value = 6503
$ws.Range("H306").Value = 70737
$ws.Range("H313").Value = 72991
$ws.Range("H317").Value = 61554
$ws.Range("I317").Value = 2143
$ws.Range("H318").Value = 24652
$ws.Range("I318").Value = 906
$ws.Range("H319").Value = 56096
$ws.Range("I319").Value = 1785
$ws.Range("H320").Value = 86149
$ws.Range("I320").Value = 3893
$ws.Range("H321").Value = 90179
$ws.Range("I321").Value = 2796
$ws.Range("H322").Value = 104385
$ws.Range("I322").Value = 2287
$ws.Range("H323").Value = 148033
$ws.Range("I323").Value = 2317
$ws.Range("H324").Value = 230592
$ws.Range("I324").Value = 2665
$ws.Range("H325").Value = 664115
$ws.Range("I325").Value = 5455
$ws.Range("H326").Value = 391033
$ws.Range("I326").Value = 3424
$ws.Range("H327").Value = 261513
$ws.Range("I327").Value = 3782

$wb.Save()
